# EPBDS: updated demo application.
#
# Applies the Chartis2AcordModelMappingExample.xlsx edits:
#  - Sheet1: refresh convertIntToC4/convertStringToC25 sample code &
#    signatures (int -> Integer boxed params, XmlBeanFactory instead of the
#    generated XmlBean .Factory), and move the active selection.
#  - Environment: add a new "org.openl.mapper.demo" import row, extending
#    the merged "import" label down one row, and move the active selection.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1 : update the Java snippets / method signatures shown in the demo
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

# C13: body of convertStringToC25(String, C25) - swap the generated
# XmlBean factory call for the XmlBeanFactory helper.
$ws1.Range("C13").Value = "if (destination==null) {`n   destination = (C25)XmlBeanFactory.newInstance(C25.class);`n}`ndestination.setId(source);`nreturn destination;"

# C19: signature of convertIntToC4(int, C4) -> now takes a boxed Integer.
$ws1.Range("C19").Value = "Method C4 convertIntToC4(Integer source, C4 destination)"

# C20: body of convertIntToC4(Integer, C4) - swap the generated XmlBean
# factory call for the XmlBeanFactory helper.
$ws1.Range("C20").Value = "if (destination==null) {`n   destination = (C4)XmlBeanFactory.newInstance(C4.class);`n}`ndestination.setId(String.valueOf(source));`nreturn destination;"

# C22: signature of convertIntToC4(C4, int) -> now returns a boxed Integer.
$ws1.Range("C22").Value = "Method int convertIntToC4(C4 source, Integer destination)"

# Move the active selection on Sheet1 (also clears the stale topLeftCell).
$ws1.Range("C22:D22").Select()

# ---------------------------------------------------------------------
# Environment : add a new import row (org.openl.mapper.demo)
# ---------------------------------------------------------------------
$wsEnv = $wb.Worksheets.Item("Environment")

# Push the formatting of the last "import" row (C7/D7) down to the new
# row 8, then promote C7 to the "middle of merge" style (matching C6).
$wsEnv.Range("C7").Copy()
$wsEnv.Range("C8").PasteSpecial(-4122)
$wsEnv.Range("D7").Copy()
$wsEnv.Range("D8").PasteSpecial(-4122)
$wsEnv.Range("C6").Copy()
$wsEnv.Range("C7").PasteSpecial(-4122)

# New import value.
$wsEnv.Range("D8").Value = "org.openl.mapper.demo"

# Extend the merged "import" label cell from C5:C7 down to C5:C8.
$wsEnv.Range("C5:C7").UnMerge()
$wsEnv.Range("C5:C8").Merge()

# Move the active selection on the Environment sheet.
$wsEnv.Range("D14").Select()

# Leave Sheet1 as the active/selected tab (matches the saved workbook).
$ws1.Select()
$ws1.Range("C22:D22").Select()
